$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 32: deposit date cleared, deposit amount zeroed ---
$ws.Range("D32").Value = ""
$ws.Range("E32").Value = 0

# --- Append new rows 33-36 with the latest collection/deposit entries ---
# Columns A,B,C,D,F,G hold text; column E holds the numeric deposit amount.
$ws.Range("A33:D36").NumberFormat = "@"
$ws.Range("F33:G36").NumberFormat = "@"

$ws.Range("A33").Value = "08-12-2025"
$ws.Range("B33").Value = "010965012-Medha Sub Division Office Coll."
$ws.Range("C33").Value = "Cash"
$ws.Range("D33").Value = "2025-12-08"
$ws.Range("E33").Value = 44300
$ws.Range("F33").Value = ""
$ws.Range("G33").Value = "2025-12-23"

$ws.Range("A34").Value = "08-12-2025"
$ws.Range("B34").Value = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
$ws.Range("C34").Value = "Cash"
$ws.Range("D34").Value = "2025-12-08"
$ws.Range("E34").Value = 64020
$ws.Range("F34").Value = ""
$ws.Range("G34").Value = "2025-12-23"

$ws.Range("A35").Value = "08-12-2025"
$ws.Range("B35").Value = "020965018-Kai Lalsingrao Shinde Gr.Big.Sheti Sah.Pat.Ltd. Br. Medha"
$ws.Range("C35").Value = "Cash"
$ws.Range("D35").Value = "2025-12-08"
$ws.Range("E35").Value = 54910
$ws.Range("F35").Value = ""
$ws.Range("G35").Value = "2025-12-23"

$ws.Range("A36").Value = "08-12-2025"
$ws.Range("B36").Value = "020965019-SHRI DATTATRAY MAHARAJ KALAMBE SAH. PAT. LTD.DAPAWADI"
$ws.Range("C36").Value = "Cash"
$ws.Range("D36").Value = "2025-12-08"
$ws.Range("E36").Value = 16290
$ws.Range("F36").Value = ""
$ws.Range("G36").Value = "2025-12-23"
